# Update scraped event stats ("想去人数" / "最低票价") in the
# "展览" and "全部类型" worksheets, mirroring the data refresh captured
# in the upstream gh-pages output commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1083
$ws1.Range("F4").Value = 185
$ws1.Range("G4").Value = 90
$ws1.Range("F5").Value = 3248
$ws1.Range("F7").Value = 297
$ws1.Range("F10").Value = 16
$ws1.Range("F12").Value = 110
$ws1.Range("F13").Value = 189
$ws1.Range("F14").Value = 31
$ws1.Range("F15").Value = 85
$ws1.Range("F16").Value = 2804
$ws1.Range("F17").Value = 1104
$ws1.Range("F18").Value = 6

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1083
$ws4.Range("F5").Value = 185
$ws4.Range("G5").Value = 90
$ws4.Range("F6").Value = 3248
$ws4.Range("F8").Value = 297
$ws4.Range("F12").Value = 16
$ws4.Range("F14").Value = 110
$ws4.Range("F15").Value = 189
$ws4.Range("F16").Value = 31
$ws4.Range("F17").Value = 85
$ws4.Range("F18").Value = 2804
$ws4.Range("F19").Value = 1104
$ws4.Range("F20").Value = 6

$wb.Save()
